# Radiant 2025.1 help-assignments workbook update
# - fills blank "Writer" (column D) cells with "no writer"
# - fills blank "Sub-sections" (column C) cells with "(no subsection)"
# - re-applies AutoFilter across the full table (A1:M231) then turns the
#   filter buttons back off, leaving the _FilterDatabase name pointing at
#   the whole table instead of the old C2:C7 sliver
# - scrolls/selects near the bottom of the table, matching the saved view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025.1")

# --- Column D ("Writer") : blank cells get "no writer" ------------------
$dRows = @(5,6,7,10,11,12,13,14,17,21,23,26,33,34,35,45,56,58,59,61,65,66,68,75,76,77,78,79,80,81,82,83,84,86,89,90,91,92,93,94,95,96,98,99,100,101,105,106,107,110,111,123,125,126,127,128,130,141,162,163,164,165,166,190,198,204,205,206,207,222,229)
foreach ($r in $dRows) {
    $ws.Cells.Item($r, 4).Value = "no writer"
}

# --- Column C ("Sub-sections") : blank cells get "(no subsection)" ------
$cRows = @(184,185,186,187,188,189,190,191,192,193,195,196,197,198,199,200,201,202,204,205,206,212,213,214,215,220,221,222,223,224,225,226,227,228,229,230,231)
foreach ($r in $cRows) {
    $ws.Cells.Item($r, 3).Value = "(no subsection)"
}

# --- Re-point the filter database at the whole table, filter off --------
$ws.Range("A1:M231").AutoFilter() | Out-Null
$ws.AutoFilterMode = $false
foreach ($n in $ws.Names) {
    if ($n.Name -eq "2025.1!_FilterDatabase") {
        $n.RefersTo = "='2025.1'!`$A`$1:`$M`$231"
    }
}

# --- Restore the saved selection/scroll position -------------------------
$ws.Activate()
$ws.Range("C237").Select()
